# Edit "D suite.xlsx" - Test Cases sheet
# Replace the placeholder "TBD" JIRA IDs in column B (rows 26-29) with the
# actual JIRA ticket numbers, and update the selected cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Column B (JIRA ID) for the type-ahead related test cases currently holds
# the placeholder value "TBD". Replace with the real JIRA IDs.
$ws.Range("B26").Value = "OPQA-512"
$ws.Range("B27").Value = "OPQA-516"
$ws.Range("B28").Value = "OPQA-517"
$ws.Range("B29").Value = "OPQA-518"

# Update the active selection on the sheet to reflect where editing ended.
$ws.Activate()
$ws.Range("B32").Select()
